$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 10:20"

$ws.Range("B4").Value = 164359
$ws.Range("C4").Value = 515
$ws.Range("E4").Value = 155679
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 3173

$ws.Range("B17").Value = 9772
$ws.Range("C17").Value = 154
$ws.Range("E17").Value = 9028

$ws.Range("B23").Value = 4465
$ws.Range("C23").Value = 20
$ws.Range("E23").Value = 4421

$ws.Range("A27").Value = "Dinamarca"
$ws.Range("B27").Value = 2815
$ws.Range("C27").Value = 238
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 2737
$ws.Range("F27").Value = 137
$ws.Range("H27").Value = 77

$ws.Range("A28").Value = "Malasia"
$ws.Range("B28").Value = 2626
$ws.Range("D28").Value = 479
$ws.Range("E28").Value = 2110
$ws.Range("F28").Value = 94
$ws.Range("H28").Value = 37

$ws.Range("B68").Value = 530
$ws.Range("C68").Value = 15
$ws.Range("E68").Value = 231

$ws.Range("A72").Value = "Bosnia y Herzegovina"
$ws.Range("B72").Value = 411
$ws.Range("C72").Value = 43
$ws.Range("D72").Value = 17
$ws.Range("E72").Value = 384
$ws.Range("F72").Value = 1
$ws.Range("H72").Value = 10

$ws.Range("A73").Value = "Letonia"
$ws.Range("B73").Value = 398
$ws.Range("C73").Value = 22
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 397
$ws.Range("F73").Value = 3
$ws.Range("H73").Value = 0

$ws.Range("A74").Value = "Bulgaria"
$ws.Range("C74").Value = 20
$ws.Range("E74").Value = 354
$ws.Range("F74").Value = 13
$ws.Range("H74").Value = 8

$ws.Range("A76").Value = "Eslovaquia"
$ws.Range("B76").Value = 363
$ws.Range("C76").Value = 27
$ws.Range("D76").Value = 7
$ws.Range("E76").Value = 356
$ws.Range("F76").Value = 1
$ws.Range("H76").Value = 0

$ws.Range("A77").Value = "Tunez"
$ws.Range("B77").Value = 362
$ws.Range("D77").Value = 3
$ws.Range("E77").Value = 350
$ws.Range("F77").Value = 10
$ws.Range("H77").Value = 9

$ws.Range("B95").Value = 174
$ws.Range("C95").Value = 4
$ws.Range("E95").Value = 168

$ws.Range("D99").Value = 28
$ws.Range("E99").Value = 134

$ws.Range("A111").Value = "Camboya"
$ws.Range("B111").Value = 109
$ws.Range("C111").Value = 2
$ws.Range("D111").Value = 23
$ws.Range("E111").Value = 86
$ws.Range("F111").Value = 1

$ws.Range("A112").Value = "Georgia"
$ws.Range("B112").Value = 108
$ws.Range("C112").Value = 5
$ws.Range("D112").Value = 21
$ws.Range("E112").Value = 87
$ws.Range("F112").Value = 6

$ws.Range("A113").Value = "Kirguistan"
$ws.Range("C113").Value = 13
$ws.Range("D113").Value = 3
$ws.Range("E113").Value = 104
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 0

$ws.Range("A114").Value = "Bolivia"
$ws.Range("C114").Value = 10
$ws.Range("D114").Value = 0
$ws.Range("E114").Value = 101
$ws.Range("F114").Value = 3
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 6

$ws.Range("A179").Value = "Republica del Chad"
$ws.Range("C179").Value = 2

$ws.Range("A180").Value = "San Cristobal y Nieves"

$ws.Range("A181").Value = "Antigua y Barbuda"
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 7
$ws.Range("H181").Value = 0

$ws.Range("A182").Value = "Angola"
$ws.Range("B182").Value = 7
$ws.Range("D182").Value = 1
$ws.Range("E182").Value = 4
$ws.Range("H182").Value = 2

$ws.Range("A183").Value = "Santa Sede"

$ws.Range("A184").Value = "San Martin (Parte Holandesa)"

$ws.Range("A185").Value = "Benin"
$ws.Range("E185").Value = 6
$ws.Range("H185").Value = 0

$ws.Range("A186").Value = "Cabo Verde"
$ws.Range("D186").Value = 0
$ws.Range("H186").Value = 1

$ws.Range("A187").Value = "San Bartolome"
$ws.Range("D187").Value = 1
$ws.Range("E187").Value = 5
$ws.Range("H187").Value = 0

$ws.Range("A188").Value = "Sudan"
$ws.Range("D188").Value = 0
$ws.Range("E188").Value = 4
$ws.Range("H188").Value = 2

$ws.Range("A189").Value = "Mauritania"
$ws.Range("B189").Value = 6
$ws.Range("D189").Value = 2
$ws.Range("E189").Value = 3
$ws.Range("H189").Value = 1

$ws.Range("A190").Value = "Montserrat"

$ws.Range("A191").Value = "Fiyi"

$ws.Range("A192").Value = "Islas Turcas y Caicos"

$ws.Range("A195").Value = "Gambia"

$ws.Range("A196").Value = "Nicaragua"

$ws.Range("A197").Value = "Belice"

$ws.Range("A200").Value = "Republica de Africa Central"
